$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.616.83'
$ws.Range("E2").Value = '  -1.10%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.277.76'
$ws.Range("E3").Value = '  +0.84%  '

# Row 4
$ws.Range("E4").Value = '  +0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.45'
$ws.Range("E5").Value = '  +0.81%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.642'
$ws.Range("E6").Value = '  +2.48%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.44'
$ws.Range("E7").Value = '  +1.80%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.439'
$ws.Range("E9").Value = '  -0.51%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0959'
$ws.Range("E10").Value = '  -5.90%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.87'
$ws.Range("E11").Value = '  -0.69%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.45'
$ws.Range("E12").Value = '  +1.93%  '

# Row 13
$ws.Range("E13").Value = '  -1.55%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.610.51'
$ws.Range("E14").Value = '  +0.73%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.10'
$ws.Range("E15").Value = '  -2.96%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.08'
$ws.Range("E16").Value = '  -0.58%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.824'
$ws.Range("E17").Value = '  -2.22%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.270.78'
$ws.Range("E18").Value = '  +0.08%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.519.77'
$ws.Range("E19").Value = '  -0.99%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0970'
$ws.Range("E20").Value = '  -3.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.03'
$ws.Range("E21").Value = '  -0.12%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.11'
$ws.Range("E22").Value = '  +1.79%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.57'
$ws.Range("E23").Value = '  -1.47%  '

# Row 24
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  +0.04%  '

# Row 25
$ws.Range("B25").Value = 'WEMIXToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.74'
$ws.Range("E25").Value = '  +15.11%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.43'
$ws.Range("E26").Value = '  +0.41%  '

# Row 27
$ws.Range("E27").Value = '  -2.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.72'
$ws.Range("E28").Value = '  -2.69%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.60'
$ws.Range("E29").Value = '  +1.03%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.77'
$ws.Range("E30").Value = '  +5.20%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("E31").Value = '  +3.95%  '

# Row 32
$ws.Range("E32").Value = '  -4.18%  '

# Row 33
$ws.Range("E33").Value = '  +1.04%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.96'
$ws.Range("E34").Value = '  +5.19%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0678'
$ws.Range("E35").Value = '  -0.58%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("E36").Value = '  +2.02%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.46'
$ws.Range("E37").Value = '  -1.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.59'
$ws.Range("E38").Value = '  -5.29%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.29'
$ws.Range("E39").Value = '  -0.28%  '

# Row 40
$ws.Range("E40").Value = '  -2.55%  '

# Row 41
$ws.Range("E41").Value = '  +0.13%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.75'
$ws.Range("E42").Value = '  +6.77%  '

# Row 43
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.20'
$ws.Range("E43").Value = '  -0.56%  '

# Row 44
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.44'
$ws.Range("E44").Value = '  +2.63%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.79'
$ws.Range("E45").Value = '  -0.39%  '

# Row 46
$ws.Range("B46").Value = 'Celestia'
$ws.Range("C46").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.19'
$ws.Range("E46").Value = '  +6.91%  '

# Row 47
$ws.Range("B47").Value = 'TrustWalletToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.19'
$ws.Range("E47").Value = '  -0.09%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0942'
$ws.Range("E48").Value = '  -2.22%  '

# Row 49
$ws.Range("E49").Value = '  -1.97%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.428.38'
$ws.Range("E50").Value = '  -0.35%  '

# Row 51
$ws.Range("E51").Value = '  -0.50%  '
